$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.5742780612326
$ws.Range("C2").Value = 7.796629849130533
$ws.Range("D2").Value = 15.1723429182464
$ws.Range("E2").Value = 16.60995602762954
$ws.Range("G2").Value = 55.09818725783304
$ws.Range("H2").Value = 20.72013471696813
$ws.Range("J2").Value = 9.502288842922749
$ws.Range("K2").Value = 15.08719614560593
$ws.Range("N2").Value = 22.37534714706104
$ws.Range("B3").Value = 15.37699206811062
$ws.Range("C3").Value = 7.656074183994916
$ws.Range("D3").Value = 15.11874594320919
$ws.Range("E3").Value = 16.55716943134362
$ws.Range("G3").Value = 54.94265393510328
$ws.Range("H3").Value = 20.7409365053499
$ws.Range("J3").Value = 9.51292979004193
$ws.Range("K3").Value = 14.96069891483657
$ws.Range("N3").Value = 22.42821639222286
$ws.Range("B4").Value = 15.258983038976
$ws.Range("C4").Value = 7.570858047596745
$ws.Range("D4").Value = 15.08928247495458
$ws.Range("E4").Value = 16.5285553120065
$ws.Range("G4").Value = 54.86026843848087
$ws.Range("H4").Value = 20.75719596397233
$ws.Range("J4").Value = 9.520939294784489
$ws.Range("K4").Value = 14.88642758747457
$ws.Range("N4").Value = 22.46260862216643
$ws.Range("B5").Value = 15.21173836574301
$ws.Range("C5").Value = 7.536455454470543
$ws.Range("D5").Value = 15.0781500998124
$ws.Range("E5").Value = 16.51785692855492
$ws.Range("G5").Value = 54.83000811124321
$ws.Range("H5").Value = 20.76469748064765
$ws.Range("J5").Value = 9.524574266383265
$ws.Range("K5").Value = 14.85704604179145
$ws.Range("N5").Value = 22.47710947221011
$ws.Range("B6").Value = 15.20394614891062
$ws.Range("C6").Value = 7.530763957925006
$ws.Range("D6").Value = 15.07635458829554
$ws.Range("E6").Value = 16.51613877061864
$ws.Range("G6").Value = 54.82518384044074
$ws.Range("H6").Value = 20.76599595014476
$ws.Range("J6").Value = 9.525200255935278
$ws.Range("K6").Value = 14.85222153151754
$ws.Range("N6").Value = 22.47954667516238
$ws.Range("B7").Value = 15.25834238286041
$ws.Range("C7").Value = 7.570392704772205
$ws.Range("D7").Value = 15.08912879050627
$ws.Range("E7").Value = 16.52840712551618
$ws.Range("G7").Value = 54.85984690810692
$ws.Range("H7").Value = 20.75729358836157
$ws.Range("J7").Value = 9.52098681516728
$ws.Range("K7").Value = 14.88602771749339
$ws.Range("N7").Value = 22.4628022185504
$ws.Range("B8").Value = 15.50564112140222
$ws.Range("C8").Value = 7.747971348102823
$ws.Range("D8").Value = 15.1531526504807
$ws.Range("E8").Value = 16.59097183196092
$ws.Range("G8").Value = 55.04184803366145
$ws.Range("H8").Value = 20.72658276550481
$ws.Range("J8").Value = 9.505651472887708
$ws.Range("K8").Value = 15.04289308820971
$ws.Range("N8").Value = 22.39317603711145
$ws.Range("B9").Value = 16.01272399511839
$ws.Range("C9").Value = 8.102598309397107
$ws.Range("D9").Value = 15.30564749849555
$ws.Range("E9").Value = 16.74342033091819
$ws.Range("G9").Value = 55.50196805390709
$ws.Range("H9").Value = 20.69407488458516
$ws.Range("J9").Value = 9.487292441668975
$ws.Range("K9").Value = 15.3760724178588
$ws.Range("N9").Value = 22.27193926126039
$ws.Range("B10").Value = 16.39513357430012
$ws.Range("C10").Value = 8.364120894143227
$ws.Range("D10").Value = 15.43352069489185
$ws.Range("E10").Value = 16.87300481261131
$ws.Range("G10").Value = 55.90160719712778
$ws.Range("H10").Value = 20.68714380372064
$ws.Range("J10").Value = 9.480946244701354
$ws.Range("K10").Value = 15.63447627922256
$ws.Range("N10").Value = 22.19217102267281
$ws.Range("B11").Value = 16.57046537433436
$ws.Range("C11").Value = 8.482712914542272
$ws.Range("D11").Value = 15.4949863190103
$ws.Range("E11").Value = 16.93563289384642
$ws.Range("G11").Value = 56.09642905408643
$ws.Range("H11").Value = 20.68767934397323
$ws.Range("J11").Value = 9.479608898774154
$ws.Range("K11").Value = 15.7545467946131
$ws.Range("N11").Value = 22.15789828106279
$ws.Range("B12").Value = 16.6369920497815
$ws.Range("C12").Value = 8.527519503458459
$ws.Range("D12").Value = 15.51872194459657
$ws.Range("E12").Value = 16.95986422070974
$ws.Range("G12").Value = 56.17204011383385
$ws.Range("H12").Value = 20.68841269088174
$ws.Range("J12").Value = 9.479325046350873
$ws.Range("K12").Value = 15.80033892243987
$ws.Range("N12").Value = 22.14520949127733
$ws.Range("B13").Value = 16.62265947583001
$ws.Range("C13").Value = 8.517874857388273
$ws.Range("D13").Value = 15.51358982828012
$ws.Range("E13").Value = 16.95462287392048
$ws.Range("G13").Value = 56.15567486801417
$ws.Range("H13").Value = 20.68823115550937
$ws.Range("J13").Value = 9.479376285065898
$ws.Range("K13").Value = 15.79046297414639
$ws.Range("N13").Value = 22.14792937237402
$ws.Range("B14").Value = 16.57593627966236
$ws.Range("C14").Value = 8.486401481270541
$ws.Range("D14").Value = 15.49692992837071
$ws.Range("E14").Value = 16.93761617108576
$ws.Range("G14").Value = 56.10261304122538
$ws.Range("H14").Value = 20.6877290439261
$ws.Range("J14").Value = 9.47958108715603
$ws.Range("K14").Value = 15.75830786841462
$ws.Range("N14").Value = 22.15684856403656
$ws.Range("B15").Value = 16.54733229857118
$ws.Range("C15").Value = 8.467108518720476
$ws.Range("D15").Value = 15.48678472682757
$ws.Range("E15").Value = 16.92726577825337
$ws.Range("G15").Value = 56.07034913360737
$ws.Range("H15").Value = 20.68749057930671
$ws.Range("J15").Value = 9.479735510742735
$ws.Range("K15").Value = 15.73865298632363
$ws.Range("N15").Value = 22.162349532632
$ws.Range("B16").Value = 16.38369725676143
$ws.Range("C16").Value = 8.356358946172241
$ws.Range("D16").Value = 15.42956891063136
$ws.Range("E16").Value = 16.86898481206697
$ws.Range("G16").Value = 55.88913415338388
$ws.Range("H16").Value = 20.68718304542482
$ws.Range("J16").Value = 9.481064804825667
$ws.Range("K16").Value = 15.62667683886789
$ws.Range("N16").Value = 22.19445135012683
$ws.Range("B17").Value = 16.28361726923495
$ws.Range("C17").Value = 8.288286864779957
$ws.Range("D17").Value = 15.3953032404269
$ws.Range("E17").Value = 16.83416428251633
$ws.Range("G17").Value = 55.78127539673948
$ws.Range("H17").Value = 20.68793928394325
$ws.Range("J17").Value = 9.482277010770721
$ws.Range("K17").Value = 15.55860158309732
$ws.Range("N17").Value = 22.2146605960498
$ws.Range("B18").Value = 16.22618685407977
$ws.Range("C18").Value = 8.249100939368043
$ws.Range("D18").Value = 15.37590560747123
$ws.Range("E18").Value = 16.81448363901096
$ws.Range("G18").Value = 55.72046617526494
$ws.Range("H18").Value = 20.68872144078069
$ws.Range("J18").Value = 9.483120113023984
$ws.Range("K18").Value = 15.51968615235567
$ws.Range("N18").Value = 22.22647401259802
$ws.Range("B19").Value = 16.20676673197039
$ws.Range("C19").Value = 8.235829132785909
$ws.Range("D19").Value = 15.36939174131663
$ws.Range("E19").Value = 16.80788014571628
$ws.Range("G19").Value = 55.70008923779846
$ws.Range("H19").Value = 20.68904588708499
$ws.Range("J19").Value = 9.483430633268451
$ws.Range("K19").Value = 15.50655238466483
$ws.Range("N19").Value = 22.23050639831677
$ws.Range("B20").Value = 16.29425766924354
$ws.Range("C20").Value = 8.295536986235428
$ws.Range("D20").Value = 15.39891878304698
$ws.Range("E20").Value = 16.83783514923485
$ws.Range("G20").Value = 55.792630278058
$ws.Range("H20").Value = 20.68782284669
$ws.Range("J20").Value = 9.482132873635313
$ws.Range("K20").Value = 15.56582378479964
$ws.Range("N20").Value = 22.21248966602608
$ws.Range("B21").Value = 16.58965693527308
$ws.Range("C21").Value = 8.49564910147539
$ws.Range("D21").Value = 15.50181097771939
$ws.Range("E21").Value = 16.94259757657759
$ws.Range("G21").Value = 56.11814905810958
$ws.Range("H21").Value = 20.68786212739075
$ws.Range("J21").Value = 9.479514893734962
$ws.Range("K21").Value = 15.7677441171208
$ws.Range("N21").Value = 22.15422092343232
$ws.Range("B22").Value = 16.78345952998421
$ws.Range("C22").Value = 8.625821040596797
$ws.Range("D22").Value = 15.57173075951832
$ws.Range("E22").Value = 17.01406326318291
$ws.Range("G22").Value = 56.34157828205697
$ws.Range("H22").Value = 20.69098013834259
$ws.Range("J22").Value = 9.47910103233626
$ws.Range("K22").Value = 15.90158135481466
$ws.Range("N22").Value = 22.11782664762177
$ws.Range("B23").Value = 16.67997666611269
$ws.Range("C23").Value = 8.556417052209778
$ws.Range("D23").Value = 15.53417341954695
$ws.Range("E23").Value = 16.97565115983131
$ws.Range("G23").Value = 56.2213652294357
$ws.Range("H23").Value = 20.68903306492181
$ws.Range("J23").Value = 9.479203335261273
$ws.Range("K23").Value = 15.82999154628515
$ws.Range("N23").Value = 22.13709655495867
$ws.Range("B24").Value = 16.28944680644726
$ws.Range("C24").Value = 8.292259359604191
$ws.Range("D24").Value = 15.39728325401721
$ws.Range("E24").Value = 16.8361744961186
$ws.Range("G24").Value = 55.78749299489123
$ws.Range("H24").Value = 20.68787440588416
$ws.Range("J24").Value = 9.482197582660275
$ws.Range("K24").Value = 15.56255793316039
$ws.Range("N24").Value = 22.213470536507
$ws.Range("B25").Value = 15.87354919105469
$ws.Range("C25").Value = 8.0062835068351
$ws.Range("D25").Value = 15.26156612842366
$ws.Range("E25").Value = 16.69904687893514
$ws.Range("G25").Value = 55.36656685488924
$ws.Range("H25").Value = 20.69989525433757
$ws.Range("J25").Value = 9.491004595837596
$ws.Range("K25").Value = 15.3760724178588
$ws.Range("N25").Value = 22.30310150506348
